$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted right after the header/first
# data row (i.e. at sheet row 8). Excel's own Rows().Insert() shifts the
# existing row 8 (and everything below it) down by one -- exactly how
# this would happen interactively -- growing the used range from
# A1:R71 to A1:R72.
$ws.Rows("8").Insert()

# After the insert, row 8 is blank and the data that used to live in row 8
# now lives in row 9. Populate the new row 8 with that same record (same
# market/category/quality/etc.) since it is the same weekly series, then
# overwrite just the date with the new week's value.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(8, $col).Value2 = $ws.Cells.Item(9, $col).Value2
}

# New reporting date for the freshly inserted week.
$ws.Cells.Item(8, 4).Value2 = 45111
